$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

Replace-Text "caseQuery(M2DocEvaluator.java:543)" "caseQuery(M2DocEvaluator.java:555)"
Replace-Text "doSwitch(M2DocEvaluator.java:1084)" "doSwitch(M2DocEvaluator.java:1096)"
Replace-Text "caseBlock(M2DocEvaluator.java:1300)" "caseBlock(M2DocEvaluator.java:1305)"
Replace-Text "caseDocumentTemplate(M2DocEvaluator.java:278)" "caseDocumentTemplate(M2DocEvaluator.java:283)"
Replace-Text "generate(M2DocEvaluator.java:267)" "generate(M2DocEvaluator.java:272)"
Replace-Text "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:475)" "prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:479)"
Replace-Text "generation(AbstractTemplatesTestSuite.java:384)" "generation(AbstractTemplatesTestSuite.java:388)"
